$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.582.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.172.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.03%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "401.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.33%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +4.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.14"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0897"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.17%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.676.72"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.13%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.06"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.69%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.169.11"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "54.446.41"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.30"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.22%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.24"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.83%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.67"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.80"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.171"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.30%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.05%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.25%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0506"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +12.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "36.69"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.19%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.30"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.63"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.87"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.89%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.10"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +10.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.292"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.27"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "131.82"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.02%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.47"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.92%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.102.02"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0520"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +13.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0340"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.43%  "
